$d = $word.ActiveDocument

$d.Content.Find.Execute("Aviation data services website developer", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Airline Consultant", 2)
